$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.5511400986501143
$ws.Range("B5").Value = -0.00954355574783676
$ws.Range("B6").Value = -0.2555446171936915
$ws.Range("B7").Value = -0.6066619314321233
$ws.Range("B8").Value = -0.2442342285147824
$ws.Range("B9").Value = 0.02800526388762647
$ws.Range("B10").Value = 0.2989632083880984
$ws.Range("B11").Value = 0.1898416505639596
$ws.Range("B12").Value = 1.064747157129269
$ws.Range("B13").Value = 0.4880328864343323
$ws.Range("B14").Value = -0.007236080745368856
$ws.Range("B15").Value = 0.02799058711211094
$ws.Range("B16").Value = 0.3
$ws.Range("B17").Value = 0.05221426964547216
$ws.Range("B18").Value = -0.3
$ws.Range("B19").Value = -0.3
$ws.Range("B20").Value = -0.4
$ws.Range("B21").Value = -0.2
$ws.Range("B22").Value = -0.4
$ws.Range("B23").Value = 0.1
$ws.Range("B24").Value = 0
$ws.Range("B25").Value = 0
$ws.Range("B26").Value = 0.2000000000000001
$ws.Range("B27").Value = 0
$ws.Range("B28").Value = 0.09999999999999998
$ws.Range("B29").Value = 0
$ws.Range("B30").Value = -0.2
$ws.Range("B31").Value = 0.09999999999999998
$ws.Range("B32").Value = 0.2999999999999999
$ws.Range("B33").Value = 0
$ws.Range("B34").Value = -0.3040828537337272
$ws.Range("B35").Value = 0
$ws.Range("B36").Value = 0.1
$ws.Range("B37").Value = 0.09999999999999998
$ws.Range("B38").Value = -0.015690968608089
$ws.Range("B39").Value = 0.05096808056002433
$ws.Range("B40").Value = 0.2151519277310885
$ws.Range("B41").Value = 0.1068236533051979
$ws.Range("B42").Value = 0.04431897453975564
$ws.Range("B43").Value = 0.09123448438175896
$ws.Range("B44").Value = -0.2716600973116948
$ws.Range("B45").Value = 0.228981402283536
$ws.Range("B46").Value = 0.2
$ws.Range("B47").Value = -0.1781284012960198
$ws.Range("B48").Value = -0.3
$ws.Range("B49").Value = 0.06505493203313417
$ws.Range("B50").Value = 0.1023597690241737
$ws.Range("B51").Value = -0.5425305662094128
$ws.Range("B52").Value = -5.4
$ws.Range("B53").Value = 2.03040224812923
$ws.Range("B54").Value = 0.4339050167294337
$ws.Range("B55").Value = 1.197090680270919
$ws.Range("B56").Value = -0.6441697788099752
$ws.Range("B57").Value = 0.4901411037136303
$ws.Range("B58").Value = -0.09098140646410988
$ws.Range("B59").Value = -0.3952916234765647
$ws.Range("B60").Value = -0.2380782088493735
$ws.Range("B61").Value = -0.02511582766690132
$ws.Range("B62").Value = 0.6876823391013496
$ws.Range("B63").Value = 0.3720617293507145
$ws.Range("B64").Value = -0.1123255314657629
$ws.Range("B65").Value = -0.1387895598915543
$ws.Range("B66").Value = 0.1213803088128225
$ws.Range("B67").Value = -0.07951810869463416
$ws.Range("B68").Value = 0.3863682696630121
$ws.Range("B69").Value = -0.006123215295980228
$ws.Range("B70").Value = 0.04342916022020096
$ws.Range("B71").Value = -0.2736421272901388
$ws.Range("B72").Value = -0.1489623566660376
$ws.Range("B73").Value = -0.1051295019602515

$ws.Range("A74:B82").EntireRow.Delete()
